$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 32.47042033333333
$ws.Cells.Item(2, 8).Value = 97.411261
$ws.Cells.Item(2, 9).Value = 0.5240295449207956
$ws.Cells.Item(2, 10).Value = 0.5240295449207955
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 117.044563
$ws.Cells.Item(2, 14).Value = 351.133689
$ws.Cells.Item(2, 15).Value = 0.3245365645427815
$ws.Cells.Item(2, 16).Value = 0.3245365645427815
$ws.Cells.Item(2, 17).Value = 3800.486158341314
$ws.Cells.Item(2, 18).Value = 34204.37542507183
$ws.Cells.Item(2, 19).Value = 0.1700667482275122
$ws.Cells.Item(2, 20).Value = 0.1700667482275121

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 32.47042033333333
$ws.Cells.Item(3, 8).Value = 97.411261
$ws.Cells.Item(3, 9).Value = 0.5240295449207956
$ws.Cells.Item(3, 10).Value = 0.5240295449207955
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 101.5800373333333
$ws.Cells.Item(3, 14).Value = 304.740112
$ws.Cells.Item(3, 15).Value = 0.281657135515876
$ws.Cells.Item(3, 16).Value = 0.281657135515876
$ws.Cells.Item(3, 17).Value = 3298.346509689025
$ws.Cells.Item(3, 18).Value = 29685.11858720123
$ws.Cells.Item(3, 19).Value = 0.1475966605480794
$ws.Cells.Item(3, 20).Value = 0.1475966605480793

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 32.47042033333333
$ws.Cells.Item(4, 8).Value = 97.411261
$ws.Cells.Item(4, 9).Value = 0.5240295449207956
$ws.Cells.Item(4, 10).Value = 0.5240295449207955
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 142.0267893333333
$ws.Cells.Item(4, 14).Value = 426.080368
$ws.Cells.Item(4, 15).Value = 0.3938062999413425
$ws.Cells.Item(4, 16).Value = 0.3938062999413425
$ws.Cells.Item(4, 17).Value = 4611.669548247116
$ws.Cells.Item(4, 18).Value = 41505.02593422405
$ws.Cells.Item(4, 19).Value = 0.206366136145204
$ws.Cells.Item(4, 20).Value = 0.206366136145204

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 22.43791
$ws.Cells.Item(5, 8).Value = 67.31372999999999
$ws.Cells.Item(5, 9).Value = 0.3621181261458191
$ws.Cells.Item(5, 10).Value = 0.362118126145819
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 117.044563
$ws.Cells.Item(5, 14).Value = 351.133689
$ws.Cells.Item(5, 15).Value = 0.3245365645427815
$ws.Cells.Item(5, 16).Value = 0.3245365645427815
$ws.Cells.Item(5, 17).Value = 2626.23537058333
$ws.Cells.Item(5, 18).Value = 23636.11833524997
$ws.Cells.Item(5, 19).Value = 0.1175205726180337
$ws.Cells.Item(5, 20).Value = 0.1175205726180337

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 22.43791
$ws.Cells.Item(6, 8).Value = 67.31372999999999
$ws.Cells.Item(6, 9).Value = 0.3621181261458191
$ws.Cells.Item(6, 10).Value = 0.362118126145819
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 101.5800373333333
$ws.Cells.Item(6, 14).Value = 304.740112
$ws.Cells.Item(6, 15).Value = 0.281657135515876
$ws.Cells.Item(6, 16).Value = 0.281657135515876
$ws.Cells.Item(6, 17).Value = 2279.243735481973
$ws.Cells.Item(6, 18).Value = 20513.19361933776
$ws.Cells.Item(6, 19).Value = 0.101993154128608
$ws.Cells.Item(6, 20).Value = 0.101993154128608

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 22.43791
$ws.Cells.Item(7, 8).Value = 67.31372999999999
$ws.Cells.Item(7, 9).Value = 0.3621181261458191
$ws.Cells.Item(7, 10).Value = 0.362118126145819
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 142.0267893333333
$ws.Cells.Item(7, 14).Value = 426.080368
$ws.Cells.Item(7, 15).Value = 0.3938062999413425
$ws.Cells.Item(7, 16).Value = 0.3938062999413425
$ws.Cells.Item(7, 17).Value = 3186.784316650293
$ws.Cells.Item(7, 18).Value = 28681.05884985264
$ws.Cells.Item(7, 19).Value = 0.1426043993991773
$ws.Cells.Item(7, 20).Value = 0.1426043993991773

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 7.054627
$ws.Cells.Item(8, 8).Value = 21.163881
$ws.Cells.Item(8, 9).Value = 0.1138523289333856
$ws.Cells.Item(8, 10).Value = 0.1138523289333855
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 117.044563
$ws.Cells.Item(8, 14).Value = 351.133689
$ws.Cells.Item(8, 15).Value = 0.3245365645427815
$ws.Cells.Item(8, 16).Value = 0.3245365645427815
$ws.Cells.Item(8, 17).Value = 825.7057343430009
$ws.Cells.Item(8, 18).Value = 7431.351609087009
$ws.Cells.Item(8, 19).Value = 0.03694924369723566
$ws.Cells.Item(8, 20).Value = 0.03694924369723566

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 7.054627
$ws.Cells.Item(9, 8).Value = 21.163881
$ws.Cells.Item(9, 9).Value = 0.1138523289333856
$ws.Cells.Item(9, 10).Value = 0.1138523289333855
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 101.5800373333333
$ws.Cells.Item(9, 14).Value = 304.740112
$ws.Cells.Item(9, 15).Value = 0.281657135515876
$ws.Cells.Item(9, 16).Value = 0.281657135515876
$ws.Cells.Item(9, 17).Value = 716.6092740327414
$ws.Cells.Item(9, 18).Value = 6449.483466294672
$ws.Cells.Item(9, 19).Value = 0.03206732083918867
$ws.Cells.Item(9, 20).Value = 0.03206732083918865

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 7.054627
$ws.Cells.Item(10, 8).Value = 21.163881
$ws.Cells.Item(10, 9).Value = 0.1138523289333856
$ws.Cells.Item(10, 10).Value = 0.1138523289333855
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 142.0267893333333
$ws.Cells.Item(10, 14).Value = 426.080368
$ws.Cells.Item(10, 15).Value = 0.3938062999413425
$ws.Cells.Item(10, 16).Value = 0.3938062999413425
$ws.Cells.Item(10, 17).Value = 1001.946022754245
$ws.Cells.Item(10, 18).Value = 9017.51420478821
$ws.Cells.Item(10, 19).Value = 0.04483576439696122
$ws.Cells.Item(10, 20).Value = 0.04483576439696121

Write-Output "Applied NATMI Adam17-Itgb1 updates"